# Added precondition script for Toggle

$wb = $excel.ActiveWorkbook

# Switch the active sheet to "Deal" (this becomes the workbook's active tab,
# and moves the "tabSelected" flag from the previously active sheet to this one)
$dealSheet = $wb.Worksheets.Item("Deal")
$dealSheet.Activate()

# New "Deal Received" value for the Stage column on the Deal precondition sheet
$dealSheet.Range("D2").Value = "Deal Received"

# Size the new column to fit its content, like the other header columns on this sheet
$dealSheet.Range("D1").ColumnWidth = 12.75

# Update the visible selection on the Deal sheet
$dealSheet.Range("L5").Select() | Out-Null
